# Auto-generated script applying scheduled market-data refresh to Lamia_Profits sheets.
# For each job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) this rewrites the cached
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H..N) with refreshed values,
# matching a scheduled Universalis price-sync run.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 11796
$ws.Range("J17").Value = 11796
$ws.Range("L17").Value = 35388
$ws.Range("N17").Value = -35724
# Row 33
$ws.Range("H33").Value = 552.8
$ws.Range("I33").Value = 534.55554
$ws.Range("J33").Value = 580.1667
$ws.Range("K33").Value = 534.55554
$ws.Range("L33").Value = 580.1667
$ws.Range("M33").Value = -305.55554
$ws.Range("N33").Value = -1038.1667
# Row 62
$ws.Range("H62").Value = 7123.3667
$ws.Range("I62").Value = 6295.0625
$ws.Range("K62").Value = 6295.0625
$ws.Range("M62").Value = -5671.0625
# Row 64
$ws.Range("H64").Value = 14500.75
$ws.Range("J64").Value = 13001
$ws.Range("L64").Value = 13001
$ws.Range("N64").Value = -13497
# Row 65
$ws.Range("H65").Value = 7123.3667
$ws.Range("I65").Value = 6295.0625
$ws.Range("K65").Value = 31475.3125
$ws.Range("M65").Value = -28355.3125
# Row 67
$ws.Range("H67").Value = 14500.75
$ws.Range("J67").Value = 13001
$ws.Range("L67").Value = 13001
$ws.Range("N67").Value = -14717
# Row 70
$ws.Range("H70").Value = 13404283
$ws.Range("I70").Value = 169866.17
$ws.Range("J70").Value = 22227228
$ws.Range("K70").Value = 509598.51
$ws.Range("L70").Value = 66681684
$ws.Range("M70").Value = -509328.51
$ws.Range("N70").Value = -66682224
# Row 73
$ws.Range("H73").Value = 13404283
$ws.Range("I73").Value = 169866.17
$ws.Range("J73").Value = 22227228
$ws.Range("K73").Value = 509598.51
$ws.Range("L73").Value = 66681684
$ws.Range("M73").Value = -508662.51
$ws.Range("N73").Value = -66683556
# Row 74
$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -24064
# Row 77
$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 125000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -120320
# Row 111
$ws.Range("H111").Value = 1362.6666
$ws.Range("I111").Value = 1362.6666
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 4087.9998
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -1020.9998
# Row 113
$ws.Range("H113").Value = 6534.294
$ws.Range("I113").Value = 3588.4285
$ws.Range("K113").Value = 3588.4285
$ws.Range("M113").Value = -334.4285
# Row 138
$ws.Range("H138").Value = 3262.08
$ws.Range("I138").Value = 2219.8572
$ws.Range("J138").Value = 3667.389
$ws.Range("K138").Value = 6659.571599999999
$ws.Range("L138").Value = 11002.167
$ws.Range("M138").Value = -1519.571599999999
$ws.Range("N138").Value = -21282.167
# Row 141
$ws.Range("H141").Value = 3643.1875
$ws.Range("I141").Value = 3592.5715
$ws.Range("J141").Value = 3997.5
$ws.Range("K141").Value = 10777.7145
$ws.Range("L141").Value = 11992.5
$ws.Range("M141").Value = -5597.7145
$ws.Range("N141").Value = -22352.5
# Cells removed entirely by the refresh (no HQ price data returned)
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("N111").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2673
$ws.Range("I45").Value = 2124.348
$ws.Range("K45").Value = 2124.348
$ws.Range("M45").Value = -1747.348
# Row 109
$ws.Range("H109").Value = 87598.5
$ws.Range("J109").Value = 87598.5
$ws.Range("L109").Value = 87598.5
$ws.Range("N109").Value = -90372.5
# Row 132
$ws.Range("H132").Value = 2865.6758
$ws.Range("I132").Value = 2170.6428
$ws.Range("J132").Value = 5028
$ws.Range("K132").Value = 6511.928400000001
$ws.Range("L132").Value = 15084
$ws.Range("M132").Value = -3981.928400000001
$ws.Range("N132").Value = -20144

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 3012.7144
$ws.Range("I94").Value = 2681.5
$ws.Range("K94").Value = 2681.5
$ws.Range("M94").Value = -2230.5
# Row 137
$ws.Range("H137").Value = 69999
$ws.Range("J137").Value = 69999
$ws.Range("L137").Value = 69999
$ws.Range("N137").Value = -80199

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 134.45
$ws.Range("I7").Value = 28.333334
$ws.Range("J7").Value = 179.92857
$ws.Range("K7").Value = 28.333334
$ws.Range("L7").Value = 179.92857
$ws.Range("M7").Value = 84.66666599999999
$ws.Range("N7").Value = -405.92857
# Row 22
$ws.Range("H22").Value = 3099.818
$ws.Range("I22").Value = 1508.5555
$ws.Range("J22").Value = 4201.4614
$ws.Range("K22").Value = 1508.5555
$ws.Range("L22").Value = 4201.4614
$ws.Range("M22").Value = -1158.5555
$ws.Range("N22").Value = -4901.4614
# Row 31
$ws.Range("H31").Value = 25808.223
$ws.Range("I31").Value = 2331.7354
$ws.Range("J31").Value = 98371.91
$ws.Range("K31").Value = 2331.7354
$ws.Range("L31").Value = 98371.91
$ws.Range("M31").Value = -2036.7354
$ws.Range("N31").Value = -98961.91
# Row 34
$ws.Range("H34").Value = 25808.223
$ws.Range("I34").Value = 2331.7354
$ws.Range("J34").Value = 98371.91
$ws.Range("K34").Value = 2331.7354
$ws.Range("L34").Value = 98371.91
$ws.Range("M34").Value = -2129.7354
$ws.Range("N34").Value = -98775.91
# Row 62
$ws.Range("H62").Value = 9613.799999999999
$ws.Range("I62").Value = 3733.1428
$ws.Range("K62").Value = 3733.1428
$ws.Range("M62").Value = -3109.1428
# Row 65
$ws.Range("H65").Value = 9613.799999999999
$ws.Range("I65").Value = 3733.1428
$ws.Range("K65").Value = 18665.714
$ws.Range("M65").Value = -15545.714

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 41547044
$ws.Range("I4").Value = 60555740
$ws.Range("K4").Value = 181667220
$ws.Range("M4").Value = -181667108
# Row 26
$ws.Range("H26").Value = 1189.1666
$ws.Range("J26").Value = 924.5
$ws.Range("L26").Value = 2773.5
$ws.Range("N26").Value = -3349.5
# Row 139
$ws.Range("H139").Value = 5858.909
$ws.Range("I139").Value = 4550
$ws.Range("J139").Value = 6349.75
$ws.Range("K139").Value = 13650
$ws.Range("L139").Value = 19049.25
$ws.Range("M139").Value = -8510
$ws.Range("N139").Value = -29329.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 15997
$ws.Range("I21").Value = 29996
$ws.Range("J21").Value = 1998
$ws.Range("K21").Value = 29996
$ws.Range("L21").Value = 1998
$ws.Range("M21").Value = -29823
$ws.Range("N21").Value = -2344
# Row 30
$ws.Range("H30").Value = 15997
$ws.Range("I30").Value = 29996
$ws.Range("J30").Value = 1998
$ws.Range("K30").Value = 29996
$ws.Range("L30").Value = 1998
$ws.Range("M30").Value = -29891
$ws.Range("N30").Value = -2208
# Row 113
$ws.Range("H113").Value = 2251.125
$ws.Range("J113").Value = 3665
$ws.Range("L113").Value = 3665
$ws.Range("N113").Value = -8005
# Row 122
$ws.Range("H122").Value = 26707.715
$ws.Range("I122").Value = 56999.5
$ws.Range("J122").Value = 14591
$ws.Range("K122").Value = 170998.5
$ws.Range("L122").Value = 43773
$ws.Range("M122").Value = -168548.5
$ws.Range("N122").Value = -48673
# Row 132
$ws.Range("H132").Value = 2888.1777
$ws.Range("I132").Value = 2536.5715
$ws.Range("J132").Value = 4118.8
$ws.Range("K132").Value = 7609.7145
$ws.Range("L132").Value = 12356.4
$ws.Range("M132").Value = -5079.7145
$ws.Range("N132").Value = -17416.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 33
$ws.Range("H33").Value = 29971.4
$ws.Range("I33").Value = 29971.4
$ws.Range("K33").Value = 29971.4
$ws.Range("M33").Value = -29681.4
# Row 68
$ws.Range("H68").Value = 6539.143
$ws.Range("I68").Value = 3302.3635
$ws.Range("K68").Value = 3302.3635
$ws.Range("M68").Value = -2553.3635
# Row 71
$ws.Range("H71").Value = 6539.143
$ws.Range("I71").Value = 3302.3635
$ws.Range("K71").Value = 16511.8175
$ws.Range("M71").Value = -12767.8175
# Row 122
$ws.Range("H122").Value = 240084.88
$ws.Range("I122").Value = 289574.34
$ws.Range("K122").Value = 868723.02
$ws.Range("M122").Value = -866273.02

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 10333.667
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 10750.375
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 10750.375
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -11998.375
# Row 65
$ws.Range("H65").Value = 10333.667
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 10750.375
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 53751.875
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -59991.875
# Row 81
$ws.Range("H81").Value = 3589.5
$ws.Range("I81").Value = 2987.9167
$ws.Range("J81").Value = 4792.6665
$ws.Range("K81").Value = 5975.8334
$ws.Range("L81").Value = 9585.333000000001
$ws.Range("M81").Value = -4914.8334
$ws.Range("N81").Value = -11707.333
# Row 84
$ws.Range("H84").Value = 3589.5
$ws.Range("I84").Value = 2987.9167
$ws.Range("J84").Value = 4792.6665
$ws.Range("K84").Value = 29879.167
$ws.Range("L84").Value = 47926.665
$ws.Range("M84").Value = -24575.167
$ws.Range("N84").Value = -58534.665
# Row 100
$ws.Range("H100").Value = 942.7778
$ws.Range("I100").Value = 760.7857
$ws.Range("J100").Value = 1579.75
$ws.Range("K100").Value = 1521.5714
$ws.Range("L100").Value = 3159.5
$ws.Range("M100").Value = -980.5714
$ws.Range("N100").Value = -4241.5
# Row 107
$ws.Range("H107").Value = 1856.909
$ws.Range("I107").Value = 1888.5
$ws.Range("J107").Value = 1819
$ws.Range("K107").Value = 5665.5
$ws.Range("L107").Value = 5457
$ws.Range("M107").Value = -3745.5
$ws.Range("N107").Value = -9297
# Row 132
$ws.Range("H132").Value = 2385.3872
$ws.Range("I132").Value = 1794.5
$ws.Range("J132").Value = 5458
$ws.Range("K132").Value = 5383.5
$ws.Range("L132").Value = 16374
$ws.Range("M132").Value = -2853.5
$ws.Range("N132").Value = -21434
